# Refreshes the cryptos.xlsx price/volume snapshot (GitHub Actions scheduled update).
#
# Every cell on this sheet is stored as plain text (inline strings like "299.67" and
# "  -1.67%  "), never as a Number/Percentage. Excel's COM layer auto-coerces a
# numeric-looking string assigned via .Value into a real Number (e.g. "299.20" would
# lose its trailing zero and become 299.2), so each write below is done with a leading
# apostrophe to force a text literal, then the style is reset to "Normal" so we don't
# leave a stray quote-prefix style behind (keeps formatting identical to the source).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'45.511.38"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -2.28%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.407.95"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +4.63%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.05%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'299.20"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -1.83%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'97.12"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -4.45%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.563"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -1.04%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  +0.03%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.512"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -2.10%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'34.96"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -4.41%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.0793"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +0.32%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'7.14"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -4.12%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.105"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +1.10%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'2.774.33"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +4.72%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'2.405.99"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +4.67%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("B16").Value = "'Polygon"
$ws.Range("B16").Style = "Normal"
$ws.Range("C16").Value = "'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("C16").Style = "Normal"
$ws.Range("D16").Value = "'0.846"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +3.80%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("B17").Value = "'Chainlink"
$ws.Range("B17").Style = "Normal"
$ws.Range("C17").Value = "'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("C17").Style = "Normal"
$ws.Range("D17").Value = "'14.27"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +3.17%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'45.488.31"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -2.37%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'12.86"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -1.40%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.0₃0952"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +0.70%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'6.23"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +3.56%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'67.16"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +1.49%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'241.67"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -3.45%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'2.81"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -3.04%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  +0.03%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'1.92"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -0.66%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'2.23"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -1.52%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'38.30"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -9.75%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'9.76"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -1.47%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'3.87"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +18.46%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'21.20"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +5.67%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'149.08"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +1.09%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  -3.93%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'5.52"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -2.27%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'0.0776"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -2.35%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("B36").Value = "'Kaspa"
$ws.Range("B36").Style = "Normal"
$ws.Range("C36").Value = "'https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("C36").Style = "Normal"
$ws.Range("D36").Value = "'0.114"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -1.15%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("B37").Value = "'ARBITRUM"
$ws.Range("B37").Style = "Normal"
$ws.Range("C37").Value = "'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("C37").Style = "Normal"
$ws.Range("D37").Value = "'1.98"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +11.63%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  -1.58%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'15.29"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -5.52%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'3.85"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -4.02%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.0301"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -0.64%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'3.25"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -3.38%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'1.942.01"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +7.11%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'  -0.07%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'90.87"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +2.89%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'1.73"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -13.59%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'8.75"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +10.63%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'15.65"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +17.65%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'102.33"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +6.77%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'  -4.02%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'2.643.38"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +4.59%  "
$ws.Range("E51").Style = "Normal"
